$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 208, pushing existing rows 208:265 down to 209:266
$ws.Rows(208).Insert()

# Populate the newly inserted row 208 with the new weekly price record
$ws.Range("A208").Value = 9
$ws.Range("B208").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C208").Value = "Metropolitana"
$ws.Range("D208").Value = 44551
$ws.Range("E208").Value = 13
$ws.Range("F208").Value = 100112044
$ws.Range("G208").Value = "Perejil"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 97
$ws.Range("K208").Value = 12000
$ws.Range("L208").Value = 13000
$ws.Range("M208").Value = 12495
$ws.Range("N208").Value = '$/docena de atados'
$ws.Range("O208").Value = "Región Metropolitana"
$ws.Range("P208").Value = 4165
$ws.Range("Q208").Value = 3
$ws.Range("R208").Value = "Hortaliza"
